$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 181,1
$arr[0,0] = 13181
$arr[1,0] = 12410
$arr[2,0] = 10979
$arr[3,0] = 9578
$arr[4,0] = 9578
$arr[5,0] = 9578
$arr[6,0] = 9263
$arr[7,0] = 9263
$arr[8,0] = 8947
$arr[9,0] = 8217
$arr[10,0] = 8217
$arr[11,0] = 8217
$arr[12,0] = 8217
$arr[13,0] = 8217
$arr[14,0] = 8146
$arr[15,0] = 8146
$arr[16,0] = 8146
$arr[17,0] = 8142
$arr[18,0] = 8142
$arr[19,0] = 8142
$arr[20,0] = 8142
$arr[21,0] = 8142
$arr[22,0] = 8142
$arr[23,0] = 8142
$arr[24,0] = 8142
$arr[25,0] = 8142
$arr[26,0] = 7598
$arr[27,0] = 7598
$arr[28,0] = 7598
$arr[29,0] = 7598
$arr[30,0] = 7598
$arr[31,0] = 7598
$arr[32,0] = 7598
$arr[33,0] = 7598
$arr[34,0] = 7598
$arr[35,0] = 7598
$arr[36,0] = 7598
$arr[37,0] = 7598
$arr[38,0] = 7598
$arr[39,0] = 7598
$arr[40,0] = 7598
$arr[41,0] = 7598
$arr[42,0] = 7598
$arr[43,0] = 7598
$arr[44,0] = 7598
$arr[45,0] = 7598
$arr[46,0] = 7598
$arr[47,0] = 7598
$arr[48,0] = 7598
$arr[49,0] = 7598
$arr[50,0] = 7598
$arr[51,0] = 7590
$arr[52,0] = 7590
$arr[53,0] = 7590
$arr[54,0] = 7590
$arr[55,0] = 7590
$arr[56,0] = 7590
$arr[57,0] = 7590
$arr[58,0] = 7590
$arr[59,0] = 7590
$arr[60,0] = 7590
$arr[61,0] = 7590
$arr[62,0] = 7590
$arr[63,0] = 7590
$arr[64,0] = 7590
$arr[65,0] = 7590
$arr[66,0] = 7590
$arr[67,0] = 7590
$arr[68,0] = 7590
$arr[69,0] = 7590
$arr[70,0] = 7573
$arr[71,0] = 7573
$arr[72,0] = 7573
$arr[73,0] = 7573
$arr[74,0] = 7573
$arr[75,0] = 7573
$arr[76,0] = 7573
$arr[77,0] = 7573
$arr[78,0] = 7573
$arr[79,0] = 7573
$arr[80,0] = 7573
$arr[81,0] = 7573
$arr[82,0] = 7573
$arr[83,0] = 7573
$arr[84,0] = 7573
$arr[85,0] = 7573
$arr[86,0] = 7573
$arr[87,0] = 7573
$arr[88,0] = 7573
$arr[89,0] = 7573
$arr[90,0] = 7573
$arr[91,0] = 7573
$arr[92,0] = 7573
$arr[93,0] = 7573
$arr[94,0] = 7573
$arr[95,0] = 7573
$arr[96,0] = 7573
$arr[97,0] = 7573
$arr[98,0] = 7573
$arr[99,0] = 7573
$arr[100,0] = 7573
$arr[101,0] = 7573
$arr[102,0] = 7573
$arr[103,0] = 7573
$arr[104,0] = 7573
$arr[105,0] = 7573
$arr[106,0] = 7573
$arr[107,0] = 7573
$arr[108,0] = 7573
$arr[109,0] = 7573
$arr[110,0] = 7573
$arr[111,0] = 7573
$arr[112,0] = 7573
$arr[113,0] = 7573
$arr[114,0] = 7573
$arr[115,0] = 7573
$arr[116,0] = 7573
$arr[117,0] = 7573
$arr[118,0] = 7573
$arr[119,0] = 7573
$arr[120,0] = 7573
$arr[121,0] = 7573
$arr[122,0] = 7573
$arr[123,0] = 7573
$arr[124,0] = 7573
$arr[125,0] = 7573
$arr[126,0] = 7573
$arr[127,0] = 7573
$arr[128,0] = 7573
$arr[129,0] = 7573
$arr[130,0] = 7573
$arr[131,0] = 7573
$arr[132,0] = 7573
$arr[133,0] = 7573
$arr[134,0] = 7573
$arr[135,0] = 7573
$arr[136,0] = 7573
$arr[137,0] = 7573
$arr[138,0] = 7573
$arr[139,0] = 7573
$arr[140,0] = 7573
$arr[141,0] = 7573
$arr[142,0] = 7573
$arr[143,0] = 7573
$arr[144,0] = 7573
$arr[145,0] = 7573
$arr[146,0] = 7573
$arr[147,0] = 7573
$arr[148,0] = 7573
$arr[149,0] = 7573
$arr[150,0] = 7573
$arr[151,0] = 7573
$arr[152,0] = 7573
$arr[153,0] = 7573
$arr[154,0] = 7573
$arr[155,0] = 7573
$arr[156,0] = 7573
$arr[157,0] = 7573
$arr[158,0] = 7573
$arr[159,0] = 7573
$arr[160,0] = 7573
$arr[161,0] = 7573
$arr[162,0] = 7573
$arr[163,0] = 7573
$arr[164,0] = 7573
$arr[165,0] = 7573
$arr[166,0] = 7573
$arr[167,0] = 7573
$arr[168,0] = 7573
$arr[169,0] = 7573
$arr[170,0] = 7573
$arr[171,0] = 7573
$arr[172,0] = 7573
$arr[173,0] = 7573
$arr[174,0] = 7573
$arr[175,0] = 7573
$arr[176,0] = 7573
$arr[177,0] = 7573
$arr[178,0] = 7573
$arr[179,0] = 7573
$arr[180,0] = 7573

$ws.Range("C2:C182").Value = $arr
